# Criação da função pred_prox_30_dias
# Updates the two forecast values produced by the new prediction
# function and moves the sheet's view/selection to the bottom of the
# data (row 15), scrolled so row 4 is at the top of the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New predicted prices (pred_prox_30_dias) for 2025-10-27 and 2025-10-29
$ws.Range("C12").Value = 360
$ws.Range("C14").Value = 369

# Scroll the window so row 4 is the top visible row, then select C15
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C15").Select()
